# Update the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions scheduled refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.814.75"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.471.82"
$ws.Range("E3").Value = "  +0.61%  "

# Row 5 - BNB
$ws.Range("D5").Value = "560.56"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6 - Solana
$ws.Range("D6").Value = "164.46"
$ws.Range("E6").Value = "  +0.34%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.96%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +5.51%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.72%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.333"

# Row 12 - Toncoin
$ws.Range("D12").Value = "4.84"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "68.726.07"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +0.99%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "23.62"
$ws.Range("E15").Value = "  +1.26%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  -2.62%  "

# Row 17 - BitcoinCash
$ws.Range("D17").Value = "338.63"
$ws.Range("E17").Value = "  -1.01%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "6.94"
$ws.Range("E18").Value = "  -2.88%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.20%  "

# Row 20 & 21 - Dai / SuiNetwork swap positions
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "1.88"
$ws.Range("E20").Value = "  +1.11%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.10%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "66.97"
$ws.Range("E22").Value = "  -0.66%  "

# Row 23 - NEARProtocol
$ws.Range("E23").Value = "  -0.42%  "

# Row 24 - Aptos
$ws.Range("D24").Value = "8.28"
$ws.Range("E24").Value = "  +2.16%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  -0.26%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "7.25"
$ws.Range("E26").Value = "  +0.79%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  -0.02%  "

# Row 28 - Bittensor
$ws.Range("D28").Value = "430.28"

# Row 29 - Fetch.AI
$ws.Range("D29").Value = "1.14"
$ws.Range("E29").Value = "  -1.40%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.76%  "

# Row 31 - Monero
$ws.Range("D31").Value = "159.81"
$ws.Range("E31").Value = "  +2.15%  "

# Row 32 - WhiteBITCoin
$ws.Range("D32").Value = "19.01"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  +0.00%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  -1.41%  "

# Row 35 - EthereumClassic
$ws.Range("D35").Value = "17.88"
$ws.Range("E35").Value = "  +0.13%  "

# Row 36 - RenderToken
$ws.Range("E36").Value = "  +0.23%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  -2.11%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -1.71%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -0.49%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  -0.09%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +1.46%  "

# Row 42 - Aave
$ws.Range("D42").Value = "130.82"
$ws.Range("E42").Value = "  -2.82%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "0.0720"
$ws.Range("E43").Value = "  +0.48%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  +1.19%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "0.566"
$ws.Range("E45").Value = "  +0.30%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  +1.46%  "

# Row 47 - BitgetToken
$ws.Range("E47").Value = "  +0.28%  "

# Row 48 - Optimism
$ws.Range("E48").Value = "  -1.82%  "

# Row 49 - THORChain
$ws.Range("E49").Value = "  -6.64%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "16.93"
$ws.Range("E50").Value = "  -2.88%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -6.28%  "
